# Fixed bug that would copy empty columns to the course equiv GUI
#
# The first data column (B) was effectively an empty/unused leading
# column that was being carried along into the generated GUI. Move its
# contents (B1:B2) out to the end of the table (G1:G2) so the visible
# range starts at column C instead of column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move (cut + paste) column B's two data cells to column G, leaving
# column B blank and shifting the used range to C1:G3.
$ws.Range("B1:B2").Cut($ws.Range("G1:G2"))
